$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L17").Value = 5061.3333
$ws.Range("H17").Value = 1640.4546
$ws.Range("N17").Value = -5397.3333
$ws.Range("J17").Value = 1687.1111
$ws.Range("I28").Value = 38462140
$ws.Range("L28").Value = 948.5
$ws.Range("M28").Value = -38461655
$ws.Range("K28").Value = 38462140
$ws.Range("J28").Value = 948.5
$ws.Range("H28").Value = 33333980
$ws.Range("N28").Value = -1918.5
$ws.Range("H38").Value = 3153.353
$ws.Range("L38").Value = 12232.8465
$ws.Range("J38").Value = 4077.6155
$ws.Range("N38").Value = -12976.8465
$ws.Range("K80").Value = 2800.875
$ws.Range("M80").Value = -1802.875
$ws.Range("H80").Value = 1065.3914
$ws.Range("I80").Value = 933.625
$ws.Range("J80").Value = 1135.6666
$ws.Range("L80").Value = 3406.9998
$ws.Range("N80").Value = -5402.9998
$ws.Range("K83").Value = 8402.625
$ws.Range("M83").Value = -3410.625
$ws.Range("H83").Value = 1065.3914
$ws.Range("L83").Value = 10220.9994
$ws.Range("N83").Value = -20204.9994
$ws.Range("I83").Value = 933.625
$ws.Range("J83").Value = 1135.6666
$ws.Range("K115").Value = 20863.5
$ws.Range("H115").Value = 6954.5
$ws.Range("I115").Value = 6954.5
$ws.Range("M115").Value = -19296.5
$ws.Range("K135").Value = 16769.5722
$ws.Range("I135").Value = 1863.2858
$ws.Range("H135").Value = 2155.1875
$ws.Range("M135").Value = -14234.5722
$ws.Range("N138").Value = -16913.2145
$ws.Range("H138").Value = 1751.17
$ws.Range("M138").Value = 1642.477
$ws.Range("J138").Value = 2211.0715
$ws.Range("I138").Value = 1165.841
$ws.Range("K138").Value = 3497.523
$ws.Range("L138").Value = 6633.2145
$ws.Range("L140").Value = 133169.8
$ws.Range("J140").Value = 133169.8
$ws.Range("H140").Value = 133169.8
$ws.Range("N140").Value = -143529.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K30").Value = 5171.6
$ws.Range("H30").Value = 5171.6
$ws.Range("M30").Value = -5021.6
$ws.Range("I30").Value = 5171.6
$ws.Range("H32").Value = 3304.68
$ws.Range("L32").Value = 10327.875
$ws.Range("M32").Value = -2406.9673
$ws.Range("K32").Value = 2693.9673
$ws.Range("I32").Value = 2693.9673
$ws.Range("N32").Value = -10901.875
$ws.Range("J32").Value = 10327.875
$ws.Range("L45").Value = 13851.167
$ws.Range("I45").Value = 8650
$ws.Range("K45").Value = 8650
$ws.Range("J45").Value = 13851.167
$ws.Range("M45").Value = -8273
$ws.Range("H45").Value = 9656.678
$ws.Range("N45").Value = -14605.167
$ws.Range("M61").Value = -2418.125
$ws.Range("I61").Value = 2630.125
$ws.Range("K61").Value = 2630.125
$ws.Range("H61").Value = 3567.2424
$ws.Range("N63").Value = -5153.7273
$ws.Range("L63").Value = 3781.7273
$ws.Range("I63").Value = 2167.2727
$ws.Range("M63").Value = -1481.2727
$ws.Range("H63").Value = 2974.5
$ws.Range("K63").Value = 2167.2727
$ws.Range("J63").Value = 3781.7273
$ws.Range("J66").Value = 3781.7273
$ws.Range("I66").Value = 2167.2727
$ws.Range("N66").Value = -25772.6365
$ws.Range("L66").Value = 18908.6365
$ws.Range("M66").Value = -7404.363499999999
$ws.Range("K66").Value = 10836.3635
$ws.Range("H66").Value = 2974.5
$ws.Range("J113").Value = 59998.75
$ws.Range("L113").Value = 59998.75
$ws.Range("H113").Value = 59998.75
$ws.Range("N113").Value = -68676.75
$ws.Range("I136").Value = 2630.125
$ws.Range("K136").Value = 7890.375
$ws.Range("M136").Value = -5340.375
$ws.Range("H136").Value = 3567.2424

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N62").ClearContents()
$ws.Range("L62").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H68").Value = 35000
$ws.Range("I68").Value = 35000
$ws.Range("K68").Value = 35000
$ws.Range("M68").Value = -34189
$ws.Range("N68").ClearContents()
$ws.Range("J68").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("J69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("I71").Value = 35000
$ws.Range("J71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("K71").Value = 105000
$ws.Range("M71").Value = -100944
$ws.Range("H71").Value = 35000
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("J72").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("L94").Value = 1719.9231
$ws.Range("J94").Value = 1719.9231
$ws.Range("N94").Value = -2621.9231
$ws.Range("H94").Value = 1106.1389

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J52").Value = 84997.5
$ws.Range("H52").Value = 84997.5
$ws.Range("L52").Value = 84997.5
$ws.Range("N52").Value = -85585.5
$ws.Range("I58").Value = 1370.5
$ws.Range("K58").Value = 1370.5
$ws.Range("H58").Value = 1690.9131
$ws.Range("M58").Value = -1167.5
$ws.Range("H62").Value = 41673184
$ws.Range("M62").Value = -62505700
$ws.Range("K62").Value = 62506324
$ws.Range("I62").Value = 62506324
$ws.Range("H65").Value = 41673184
$ws.Range("M65").Value = -312528500
$ws.Range("I65").Value = 62506324
$ws.Range("K65").Value = 312531620
$ws.Range("L107").Value = 1929.8889
$ws.Range("H107").Value = 1067.1212
$ws.Range("J107").Value = 1929.8889
$ws.Range("N107").Value = -5769.8889
$ws.Range("I136").Value = 1370.5
$ws.Range("K136").Value = 4111.5
$ws.Range("M136").Value = -1561.5
$ws.Range("H136").Value = 1690.9131

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K9").Value = 11000223.6
$ws.Range("N9").Value = -432956.74
$ws.Range("I9").Value = 3666741.2
$ws.Range("L9").Value = 432508.74
$ws.Range("H9").Value = 1200941.1
$ws.Range("J9").Value = 144169.58
$ws.Range("M9").Value = -10999999.6
$ws.Range("L68").Value = 5944.8462
$ws.Range("H68").Value = 1694.8334
$ws.Range("N68").Value = -7566.8462
$ws.Range("J68").Value = 1981.6154
$ws.Range("L71").Value = 17834.5386
$ws.Range("J71").Value = 1981.6154
$ws.Range("N71").Value = -25946.5386
$ws.Range("H71").Value = 1694.8334
$ws.Range("H81").Value = 2553.5715
$ws.Range("K81").Value = 6771.428400000001
$ws.Range("L81").Value = 8550
$ws.Range("N81").Value = -10796
$ws.Range("M81").Value = -5648.428400000001
$ws.Range("I81").Value = 2257.1428
$ws.Range("J81").Value = 2850
$ws.Range("M84").Value = -14698.2852
$ws.Range("N84").Value = -36882
$ws.Range("L84").Value = 25650
$ws.Range("H84").Value = 2553.5715
$ws.Range("J84").Value = 2850
$ws.Range("K84").Value = 20314.2852
$ws.Range("I84").Value = 2257.1428
$ws.Range("I114").Value = 1099.5
$ws.Range("J114").Value = 4450
$ws.Range("K114").Value = 3298.5
$ws.Range("H114").Value = 2216.3333
$ws.Range("L114").Value = 13350
$ws.Range("M114").Value = -44.5
$ws.Range("N114").Value = -19858
$ws.Range("I131").Value = 4902836
$ws.Range("K131").Value = 14708508
$ws.Range("M131").Value = -14703468
$ws.Range("H131").Value = 3677381
$ws.Range("N132").Value = -26793.0704
$ws.Range("I132").Value = 1228.9231
$ws.Range("H132").Value = 2038.7805
$ws.Range("J132").Value = 2414.7856
$ws.Range("L132").Value = 21733.0704
$ws.Range("M132").Value = -8530.3079
$ws.Range("K132").Value = 11060.3079

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N22").Value = -3558.75
$ws.Range("I22").Value = 1162.6666
$ws.Range("K22").Value = 1162.6666
$ws.Range("H22").Value = 2194.7144
$ws.Range("J22").Value = 2968.75
$ws.Range("L22").Value = 2968.75
$ws.Range("M22").Value = -867.6666
$ws.Range("I27").Value = 1162.6666
$ws.Range("K27").Value = 1162.6666
$ws.Range("H27").Value = 2194.7144
$ws.Range("N27").Value = -3182.75
$ws.Range("M27").Value = -1055.6666
$ws.Range("L27").Value = 2968.75
$ws.Range("J27").Value = 2968.75
$ws.Range("I93").Value = 22225286
$ws.Range("N93").Value = -16184.5
$ws.Range("K93").Value = 22225286
$ws.Range("M93").Value = -22224038
$ws.Range("J93").Value = 13688.5
$ws.Range("H93").Value = 15390949
$ws.Range("L93").Value = 13688.5
$ws.Range("I136").Value = 6671867
$ws.Range("K136").Value = 20015601
$ws.Range("M136").Value = -20013051
$ws.Range("L136").Value = 13840.3329
$ws.Range("N136").Value = -18940.3329
$ws.Range("H136").Value = 5005053.5
$ws.Range("J136").Value = 4613.4443

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M113").Value = 46
$ws.Range("I113").Value = 708
$ws.Range("H113").Value = 801.9524
$ws.Range("K113").Value = 2124
$ws.Range("J126").Value = 2469.6667
$ws.Range("H126").Value = 2479.8667
$ws.Range("N126").Value = -12349.0001
$ws.Range("K126").Value = 7447.250100000001
$ws.Range("I126").Value = 2482.4167
$ws.Range("L126").Value = 7409.000100000001
$ws.Range("M126").Value = -4977.250100000001
